$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

# Move the "Sprint 1" comment from E7 to E11
$ws.Range("E7").ClearContents()
$ws.Range("E11").Value = "Sprint 1"

# Update the selection to match the new active cell / selected range
$ws.Range("A2:E15").Select()
